$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.970.50"
$ws.Range("E2").Value = "  +5.84%  "
$ws.Range("D3").Value = "3.116.48"
$ws.Range("E3").Value = "  +3.78%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.02"
$ws.Range("E5").Value = "  +4.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.64"
$ws.Range("E6").Value = "  +3.06%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.102.54"
$ws.Range("E8").Value = "  +3.64%  "
$ws.Range("E9").Value = "  +2.39%  "
$ws.Range("E10").Value = "  +9.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.73"
$ws.Range("E11").Value = "  +9.82%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.469"
$ws.Range("E12").Value = "  +1.99%  "
$ws.Range("E13").Value = "  +5.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.56"
$ws.Range("E14").Value = "  +5.09%  "
$ws.Range("E15").Value = "  +0.72%  "
$ws.Range("D16").Value = "3.632.85"
$ws.Range("E17").Value = "  -1.04%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.114.16"
$ws.Range("E18").Value = "  +3.72%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "62.941.57"
$ws.Range("E19").Value = "  +5.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "453.09"
$ws.Range("E20").Value = "  +4.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.09"
$ws.Range("E21").Value = "  +3.10%  "
$ws.Range("E22").Value = "  +1.54%  "
$ws.Range("E23").Value = "  +5.56%  "
$ws.Range("E24").Value = "  +0.73%  "
$ws.Range("E25").Value = "  +1.86%  "
$ws.Range("E26").Value = "  +0.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.26"
$ws.Range("E27").Value = "  +0.46%  "
$ws.Range("E28").Value = "  +5.86%  "
$ws.Range("B29").Value = "FirstDigitalUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.28"
$ws.Range("E30").Value = "  +4.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.87"
$ws.Range("E31").Value = "  +11.38%  "
$ws.Range("E32").Value = "  +11.84%  "
$ws.Range("E33").Value = "  +5.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.04"
$ws.Range("E34").Value = "  +4.40%  "
$ws.Range("D35").Value = "0.0₃0807"
$ws.Range("E35").Value = "  +5.50%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.31"
$ws.Range("E36").Value = "  +8.80%  "
$ws.Range("E37").Value = "  +1.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "50.96"
$ws.Range("E38").Value = "  +3.98%  "
$ws.Range("E39").Value = "  +9.76%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.79"
$ws.Range("E40").Value = "  +1.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "428.48"
$ws.Range("E41").Value = "  +4.13%  "
$ws.Range("D42").Value = "2.962.63"
$ws.Range("E42").Value = "  +6.55%  "
$ws.Range("E43").Value = "  +5.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.111"
$ws.Range("E44").Value = "  +3.23%  "
$ws.Range("E45").Value = "  +8.80%  "
$ws.Range("E46").Value = "  +7.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "125.15"
$ws.Range("E47").Value = "  +1.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.71"
$ws.Range("E49").Value = "  -0.84%  "
$ws.Range("E50").Value = "  +0.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.80"
$ws.Range("E51").Value = "  +4.94%  "
